$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Mapping Type" labels -------------------------------------
# "PROCESSOR_EXTENSION_SOURCE" -> "EXTENSION_SOURCE"
# (shared by rows 16 and 21 - both mapping-type cells carry this text)
$ws.Range("B16").Value = "EXTENSION_SOURCE"
$ws.Range("B21").Value = "EXTENSION_SOURCE"

# "XTENSION_SOURCE_TARGET" -> "EXTENSION_SOURCE_TARGET"
$ws.Range("B25").Value = "EXTENSION_SOURCE_TARGET"

# --- C26 rich text edit: "bus_id" -> "bus_c8ySourceId" ----------------
# Edit only the affected characters, then restore per-run formatting
# (font/size/color for the whole text, strikethrough only on the
# "check: Use external id" line) since character-level writes collapse
# existing run formatting.
$full = $ws.Range("C26").Text
$idx = $full.IndexOf("bus_id")
$busIdChars = $ws.Range("C26").Characters($idx + 1, 6)
$busIdChars.Text = "bus_c8ySourceId"

$full2 = $ws.Range("C26").Text
$checkStart = $full2.IndexOf("check: Use external id")
$checkLen = "check: Use external id`n".Length

$part1 = $ws.Range("C26").Characters(1, $checkStart)
$part1.Font.Size = 9
$part1.Font.Color = 0
$part1.Font.Name = "Courier New"

$chk = $ws.Range("C26").Characters($checkStart + 1, $checkLen)
$chk.Font.Strikethrough = $true
$chk.Font.Size = 9
$chk.Font.Color = 0
$chk.Font.Name = "Courier New"

$part3Start = $checkStart + $checkLen
$part3Len = $full2.Length - $part3Start
$part3 = $ws.Range("C26").Characters($part3Start + 1, $part3Len)
$part3.Font.Size = 9
$part3.Font.Color = 0
$part3.Font.Name = "Courier New"

# --- New cell style for the "Mapping Type" header (B1) -----------------
# Same as its current style but without wrap text - this mints the new
# cellXfs entry and re-points B1 at it.
$ws.Range("B1").WrapText = $false

# --- Column width adjustments (navigation / connector column got wider) -
$ws.Columns.Item(2).ColumnWidth = 23.45   # -> stored width ~24.25 (18.875 before)
$ws.Columns.Item(3).ColumnWidth = 53.15   # -> stored width 54 (51.125 before)

# --- Scroll position / selection ---------------------------------------
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B41").Select() | Out-Null
